$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Item ID 5503) on ALC
$ws.Cells.Item(5, 8).Value = 103.22222
$ws.Cells.Item(5, 9).Value = 103.22222
$ws.Cells.Item(5, 11).Value = 103.22222
$ws.Cells.Item(5, 13).Value = 11.77778000000001

# Row 9 (Item ID 5487) on ALC
$ws.Cells.Item(9, 8).Value = 18.285715
$ws.Cells.Item(9, 9).Value = 18.285715
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 18.285715
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 150.714285
$ws.Cells.Item(9, 14).Value = ""

# Row 76 (Item ID 12602) on ALC
$ws.Cells.Item(76, 8).Value = 5939.2666
$ws.Cells.Item(76, 9).Value = 4179.8
$ws.Cells.Item(76, 11).Value = 4179.8
$ws.Cells.Item(76, 13).Value = -3864.8

# Row 79 (Item ID 12602) on ALC
$ws.Cells.Item(79, 8).Value = 5939.2666
$ws.Cells.Item(79, 9).Value = 4179.8
$ws.Cells.Item(79, 11).Value = 4179.8
$ws.Cells.Item(79, 13).Value = -3087.8

# Row 129 (Item ID 36115) on ALC
$ws.Cells.Item(129, 8).Value = 2729.75
$ws.Cells.Item(129, 9).Value = 505.7143
$ws.Cells.Item(129, 11).Value = 1517.1429
$ws.Cells.Item(129, 13).Value = 3482.8571

# Row 141 (Item ID 44161) on ALC
$ws.Cells.Item(141, 8).Value = 3428.0833
$ws.Cells.Item(141, 9).Value = 2648.818
$ws.Cells.Item(141, 11).Value = 7946.454000000001
$ws.Cells.Item(141, 13).Value = -2766.454000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Item ID 43999) on ARM
$ws.Cells.Item(61, 8).Value = 4582.1816
$ws.Cells.Item(61, 9).Value = 4536.375
$ws.Cells.Item(61, 11).Value = 4536.375
$ws.Cells.Item(61, 13).Value = -4324.375

# Row 74 (Item ID 44000) on ARM
$ws.Cells.Item(74, 8).Value = 1532.4706
$ws.Cells.Item(74, 9).Value = 1168.9181
$ws.Cells.Item(74, 10).Value = 4700.5713
$ws.Cells.Item(74, 11).Value = 1168.9181
$ws.Cells.Item(74, 12).Value = 4700.5713
$ws.Cells.Item(74, 13).Value = -294.9181000000001
$ws.Cells.Item(74, 14).Value = -6448.5713

# Row 77 (Item ID 44000) on ARM
$ws.Cells.Item(77, 8).Value = 1532.4706
$ws.Cells.Item(77, 9).Value = 1168.9181
$ws.Cells.Item(77, 10).Value = 4700.5713
$ws.Cells.Item(77, 11).Value = 5844.5905
$ws.Cells.Item(77, 12).Value = 23502.8565
$ws.Cells.Item(77, 13).Value = -1476.5905
$ws.Cells.Item(77, 14).Value = -32238.8565

# Row 97 (Item ID 19941) on ARM
$ws.Cells.Item(97, 8).Value = 1741.8889
$ws.Cells.Item(97, 10).Value = 10010
$ws.Cells.Item(97, 12).Value = 10010
$ws.Cells.Item(97, 14).Value = -11002

# Row 102 (Item ID 19945) on ARM
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = ""
$ws.Cells.Item(102, 14).Value = ""

# Row 122 (Item ID 36168) on ARM
$ws.Cells.Item(122, 8).Value = 3708.3333
$ws.Cells.Item(122, 9).Value = 2500
$ws.Cells.Item(122, 10).Value = 4312.5
$ws.Cells.Item(122, 11).Value = 7500
$ws.Cells.Item(122, 12).Value = 12937.5
$ws.Cells.Item(122, 13).Value = -5050
$ws.Cells.Item(122, 14).Value = -17837.5

# Row 123 (Item ID 34107) on ARM
$ws.Cells.Item(123, 8).Value = 16699
$ws.Cells.Item(123, 10).Value = 16699
$ws.Cells.Item(123, 12).Value = 16699
$ws.Cells.Item(123, 14).Value = -26499

# Row 125 (Item ID 34251) on ARM
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = ""

# Row 136 (Item ID 43999) on ARM
$ws.Cells.Item(136, 8).Value = 4582.1816
$ws.Cells.Item(136, 9).Value = 4536.375
$ws.Cells.Item(136, 11).Value = 13609.125
$ws.Cells.Item(136, 13).Value = -11059.125

$ws = $wb.Worksheets.Item("BSM")
# Row 80 (Item ID 13747) on BSM
$ws.Cells.Item(80, 8).Value = 285.5
$ws.Cells.Item(80, 9).Value = 339.33334
$ws.Cells.Item(80, 10).Value = 188.6
$ws.Cells.Item(80, 11).Value = 339.33334
$ws.Cells.Item(80, 12).Value = 188.6
$ws.Cells.Item(80, 13).Value = 658.66666
$ws.Cells.Item(80, 14).Value = -2184.6

# Row 83 (Item ID 13747) on BSM
$ws.Cells.Item(83, 8).Value = 285.5
$ws.Cells.Item(83, 9).Value = 339.33334
$ws.Cells.Item(83, 10).Value = 188.6
$ws.Cells.Item(83, 11).Value = 1696.6667
$ws.Cells.Item(83, 12).Value = 943
$ws.Cells.Item(83, 13).Value = 3295.3333
$ws.Cells.Item(83, 14).Value = -10927

# Row 107 (Item ID 27706) on BSM
$ws.Cells.Item(107, 8).Value = 603.3929000000001
$ws.Cells.Item(107, 9).Value = 592.44446
$ws.Cells.Item(107, 10).Value = 899
$ws.Cells.Item(107, 11).Value = 592.44446
$ws.Cells.Item(107, 12).Value = 899
$ws.Cells.Item(107, 13).Value = 1327.55554
$ws.Cells.Item(107, 14).Value = -4739

$ws = $wb.Worksheets.Item("CRP")
# Row 6 (Item ID 2219) on CRP
$ws.Cells.Item(6, 8).Value = 6000346
$ws.Cells.Item(6, 9).Value = 6666910.5
$ws.Cells.Item(6, 11).Value = 6666910.5
$ws.Cells.Item(6, 13).Value = -6666797.5

# Row 31 (Item ID 44023) on CRP
$ws.Cells.Item(31, 8).Value = 3757.8965
$ws.Cells.Item(31, 9).Value = 1635.6111
$ws.Cells.Item(31, 11).Value = 1635.6111
$ws.Cells.Item(31, 13).Value = -1340.6111

# Row 34 (Item ID 44023) on CRP
$ws.Cells.Item(34, 8).Value = 3757.8965
$ws.Cells.Item(34, 9).Value = 1635.6111
$ws.Cells.Item(34, 11).Value = 1635.6111
$ws.Cells.Item(34, 13).Value = -1433.6111

# Row 58 (Item ID 44021) on CRP
$ws.Cells.Item(58, 8).Value = 3927.8572
$ws.Cells.Item(58, 9).Value = 2064.8
$ws.Cells.Item(58, 10).Value = 5621.5454
$ws.Cells.Item(58, 11).Value = 2064.8
$ws.Cells.Item(58, 12).Value = 5621.5454
$ws.Cells.Item(58, 13).Value = -1861.8
$ws.Cells.Item(58, 14).Value = -6027.5454

# Row 62 (Item ID 12580) on CRP
$ws.Cells.Item(62, 8).Value = 49509.445
$ws.Cells.Item(62, 9).Value = 5598.5
$ws.Cells.Item(62, 10).Value = 62055.43
$ws.Cells.Item(62, 11).Value = 5598.5
$ws.Cells.Item(62, 12).Value = 62055.43
$ws.Cells.Item(62, 13).Value = -4974.5
$ws.Cells.Item(62, 14).Value = -63303.43

# Row 65 (Item ID 12580) on CRP
$ws.Cells.Item(65, 8).Value = 49509.445
$ws.Cells.Item(65, 9).Value = 5598.5
$ws.Cells.Item(65, 10).Value = 62055.43
$ws.Cells.Item(65, 11).Value = 27992.5
$ws.Cells.Item(65, 12).Value = 310277.15
$ws.Cells.Item(65, 13).Value = -24872.5
$ws.Cells.Item(65, 14).Value = -316517.15

# Row 93 (Item ID 19516) on CRP
$ws.Cells.Item(93, 8).Value = 1460.75
$ws.Cells.Item(93, 9).Value = 1460.75
$ws.Cells.Item(93, 11).Value = 1460.75
$ws.Cells.Item(93, 13).Value = 411.25

# Row 107 (Item ID 27689) on CRP
$ws.Cells.Item(107, 8).Value = 2154.9285
$ws.Cells.Item(107, 9).Value = 2427.4092
$ws.Cells.Item(107, 10).Value = 1155.8334
$ws.Cells.Item(107, 11).Value = 2427.4092
$ws.Cells.Item(107, 12).Value = 1155.8334
$ws.Cells.Item(107, 13).Value = -507.4092000000001
$ws.Cells.Item(107, 14).Value = -4995.8334

# Row 122 (Item ID 36196) on CRP
$ws.Cells.Item(122, 8).Value = 3150.3157
$ws.Cells.Item(122, 9).Value = 3255.5334
$ws.Cells.Item(122, 11).Value = 9766.600199999999
$ws.Cells.Item(122, 13).Value = -7316.600199999999

# Row 136 (Item ID 44021) on CRP
$ws.Cells.Item(136, 8).Value = 3927.8572
$ws.Cells.Item(136, 9).Value = 2064.8
$ws.Cells.Item(136, 10).Value = 5621.5454
$ws.Cells.Item(136, 11).Value = 6194.400000000001
$ws.Cells.Item(136, 12).Value = 16864.6362
$ws.Cells.Item(136, 13).Value = -3644.400000000001
$ws.Cells.Item(136, 14).Value = -21964.6362

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Item ID 4650) on CUL
$ws.Cells.Item(4, 8).Value = 139100320
$ws.Cells.Item(4, 9).Value = 139100320
$ws.Cells.Item(4, 11).Value = 417300960
$ws.Cells.Item(4, 13).Value = -417300848

# Row 38 (Item ID 4860) on CUL
$ws.Cells.Item(38, 8).Value = 93.583336
$ws.Cells.Item(38, 9).Value = 66.125
$ws.Cells.Item(38, 11).Value = 198.375
$ws.Cells.Item(38, 13).Value = 148.625

# Row 137 (Item ID 44088) on CUL
$ws.Cells.Item(137, 8).Value = 4541.0586
$ws.Cells.Item(137, 9).Value = 3248.5
$ws.Cells.Item(137, 11).Value = 9745.5
$ws.Cells.Item(137, 13).Value = -4645.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Item ID 36169) on GSM
$ws.Cells.Item(102, 8).Value = 4487
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).Value = ""

# Row 126 (Item ID 36184) on GSM
$ws.Cells.Item(126, 8).Value = 4229.2144
$ws.Cells.Item(126, 10).Value = 4512.4443
$ws.Cells.Item(126, 12).Value = 13537.3329
$ws.Cells.Item(126, 14).Value = -18477.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Item ID 36249) on LTW
$ws.Cells.Item(7, 8).Value = 2752
$ws.Cells.Item(7, 9).Value = 2752
$ws.Cells.Item(7, 11).Value = 2752
$ws.Cells.Item(7, 13).Value = -2640

# Row 46 (Item ID 5282) on LTW
$ws.Cells.Item(46, 8).Value = 2649.8572
$ws.Cells.Item(46, 9).Value = 2499.75
$ws.Cells.Item(46, 11).Value = 2499.75
$ws.Cells.Item(46, 13).Value = -2311.75

# Row 100 (Item ID 19995) on LTW
$ws.Cells.Item(100, 8).Value = 1257
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).Value = ""

# Row 122 (Item ID 36247) on LTW
$ws.Cells.Item(122, 8).Value = 24999
$ws.Cells.Item(122, 9).Value = 24999
$ws.Cells.Item(122, 11).Value = 74997
$ws.Cells.Item(122, 13).Value = -72547

# Row 126 (Item ID 36249) on LTW
$ws.Cells.Item(126, 8).Value = 2752
$ws.Cells.Item(126, 9).Value = 2752
$ws.Cells.Item(126, 11).Value = 8256
$ws.Cells.Item(126, 13).Value = -5786

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (Item ID 44031) on WVR
$ws.Cells.Item(136, 8).Value = 2355.1333
$ws.Cells.Item(136, 9).Value = 1913.3334
$ws.Cells.Item(136, 11).Value = 5740.0002
$ws.Cells.Item(136, 13).Value = -3190.0002
